$p = $ppt.ActivePresentation
$s1 = $p.Slides.Item(1)

# Duplicate slide 1 to create slide 2 (identical shapes, picture, table, group, connector)
$s1.Duplicate() | Out-Null
$s2 = $p.Slides.Item(2)

# Apply rotations to the new slide's shapes, mirroring BaseShape.rotation scenarios
$s2.Shapes.Item(1).Rotation = 10   # Rounded Rectangle 1 -> rot=600000 (10 deg)
$s2.Shapes.Item(2).Rotation = 20   # Picture 2 -> rot=1200000 (20 deg)
$s2.Shapes.Item(4).Rotation = 40   # Group 8 -> rot=2400000 (40 deg)
$s2.Shapes.Item(5).Rotation = 50   # Elbow Connector 10 -> rot=3000000 (50 deg)
